# Regenerate the "K" column (G) of the save_data sheet with recalculated
# strikeout values (previously it held a different, now-obsolete "Strike#"
# derived figure). Row 1 is the header, data rows run 2..41 for this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K), rows 2-41, per the regenerated save_data.
$newK = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 0
    6  = 3
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    11 = 0
    12 = 1
    13 = 1
    14 = 2
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 2
    27 = 1
    28 = 0
    29 = 2
    30 = 2
    31 = 0
    32 = 1
    33 = 1
    34 = 2
    35 = 2
    36 = 1
    37 = 1
    38 = 4
    39 = 1
    40 = 1
    41 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
